# Commit: "modified the Social media button (from BJ005)"
#
# Sheet1 holds finished backlog items, Sheet2 holds the remaining backlog.
# BJ-005 ("be able to forfeit a game" / Social media button story) is being
# finished, so it moves from Sheet2 row 2 down into Sheet1 as a new last row,
# and the previous last row on Sheet1 (BJ-005 "social media accounts" story)
# gets its status marked "finished" too.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# 1) Mark the current last Sheet1 row (row 6) as "finished" in column F,
#    copying the formatting already used by the rest of that row (E6).
$ws1.Range("F6").Value2 = "finished"
$ws1.Range("E6").Copy()
$ws1.Range("F6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Move the BJ-005 row from Sheet2 (A2:F2) onto Sheet1 as the new row 7,
#    carrying over values and formatting.
$src = $ws2.Range("A2:F2")
$dst = $ws1.Range("A7:F7")
$src.Copy($dst)
$excel.CutCopyMode = $false

# 3) Clear the now-relocated row on Sheet2 (contents + formatting) in place,
#    leaving the other backlog rows (and their row numbers) untouched.
$ws2.Range("A2:F2").Clear()

# 4) Restore the selections shown in each sheet, making sure Sheet1 remains
#    the active tab (select Sheet2 first, then finish on Sheet1).
$ws2.Range("A2:F2").Select()
$ws1.Activate()
$ws1.Range("C7").Select()
